$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Write row 18
$ws.Cells.Item(18, 6).Value = 'Jedinstvo'
$ws.Cells.Item(18, 7).Value = 0
$ws.Cells.Item(18, 8).Value = 'Sutjeska'
$ws.Cells.Item(18, 9).Value = 0
$ws.Cells.Item(18, 10).Value = 5.84
$ws.Cells.Item(18, 11).Value = '12/08/2023 17:12'
$ws.Cells.Item(18, 12).Value = 4.16
$ws.Cells.Item(18, 13).Value = '13/08/2023 19:20'
$ws.Cells.Item(18, 14).Value = 3.85
$ws.Cells.Item(18, 15).Value = '12/08/2023 17:12'
$ws.Cells.Item(18, 16).Value = 3.62
$ws.Cells.Item(18, 17).Value = '13/08/2023 19:20'
$ws.Cells.Item(18, 18).Value = 1.48
$ws.Cells.Item(18, 19).Value = '12/08/2023 17:12'
$ws.Cells.Item(18, 20).Value = 1.8
$ws.Cells.Item(18, 21).Value = '13/08/2023 19:20'
$ws.Cells.Item(18, 22).Value = 'https://www.betexplorer.com/football/montenegro/prva-crnogorska-liga/jedinstvo-sutjeska/YTMvEs5B/'

# Write row 19
$ws.Cells.Item(19, 6).Value = 'Mladost DG'
$ws.Cells.Item(19, 7).Value = 0
$ws.Cells.Item(19, 8).Value = 'Decic'
$ws.Cells.Item(19, 9).Value = 2
$ws.Cells.Item(19, 10).Value = 2.64
$ws.Cells.Item(19, 11).Value = '13/08/2023 10:35'
$ws.Cells.Item(19, 12).Value = 3.16
$ws.Cells.Item(19, 13).Value = '13/08/2023 19:49'
$ws.Cells.Item(19, 14).Value = 2.99
$ws.Cells.Item(19, 15).Value = '13/08/2023 10:35'
$ws.Cells.Item(19, 16).Value = 2.89
$ws.Cells.Item(19, 17).Value = '13/08/2023 19:39'
$ws.Cells.Item(19, 18).Value = 2.75
$ws.Cells.Item(19, 19).Value = '13/08/2023 10:35'
$ws.Cells.Item(19, 20).Value = 2.36
$ws.Cells.Item(19, 21).Value = '13/08/2023 19:49'
$ws.Cells.Item(19, 22).Value = 'https://www.betexplorer.com/football/montenegro/prva-crnogorska-liga/mladost-dg-decic/8nFiBukU/'

# Write row 31
$ws.Cells.Item(31, 6).Value = 'Rudar'
$ws.Cells.Item(31, 7).Value = 0
$ws.Cells.Item(31, 8).Value = 'Mornar Bar'
$ws.Cells.Item(31, 9).Value = 0
$ws.Cells.Item(31, 10).Value = 3.06
$ws.Cells.Item(31, 11).Value = '02/09/2023 06:12'
$ws.Cells.Item(31, 12).Value = 2.43
$ws.Cells.Item(31, 13).Value = '03/09/2023 16:46'
$ws.Cells.Item(31, 14).Value = 2.82
$ws.Cells.Item(31, 15).Value = '02/09/2023 06:12'
$ws.Cells.Item(31, 16).Value = 2.98
$ws.Cells.Item(31, 17).Value = '03/09/2023 16:11'
$ws.Cells.Item(31, 18).Value = 2.31
$ws.Cells.Item(31, 19).Value = '02/09/2023 06:12'
$ws.Cells.Item(31, 20).Value = 3.09
$ws.Cells.Item(31, 21).Value = '03/09/2023 16:46'
$ws.Cells.Item(31, 22).Value = 'https://www.betexplorer.com/football/montenegro/prva-crnogorska-liga/rudar-mornar-bar/8ptsWFQd/'

# Write row 32
$ws.Cells.Item(32, 6).Value = 'Arsenal Tivat'
$ws.Cells.Item(32, 7).Value = 2
$ws.Cells.Item(32, 8).Value = 'Petrovac'
$ws.Cells.Item(32, 9).Value = 2
$ws.Cells.Item(32, 10).Value = 2.53
$ws.Cells.Item(32, 11).Value = '02/09/2023 05:42'
$ws.Cells.Item(32, 12).Value = 2.66
$ws.Cells.Item(32, 13).Value = '03/09/2023 16:41'
$ws.Cells.Item(32, 14).Value = 2.84
$ws.Cells.Item(32, 15).Value = '02/09/2023 05:42'
$ws.Cells.Item(32, 16).Value = 2.87
$ws.Cells.Item(32, 17).Value = '03/09/2023 16:41'
$ws.Cells.Item(32, 18).Value = 2.73
$ws.Cells.Item(32, 19).Value = '02/09/2023 05:42'
$ws.Cells.Item(32, 20).Value = 2.88
$ws.Cells.Item(32, 21).Value = '03/09/2023 16:41'
$ws.Cells.Item(32, 22).Value = 'https://www.betexplorer.com/football/montenegro/prva-crnogorska-liga/arsenal-tivat-petrovac/f1rgTDAF/'

# Write row 61
$ws.Cells.Item(61, 6).Value = 'Sutjeska'
$ws.Cells.Item(61, 7).Value = 2
$ws.Cells.Item(61, 8).Value = 'Jedinstvo'
$ws.Cells.Item(61, 9).Value = 0
$ws.Cells.Item(61, 10).Value = 1.29
$ws.Cells.Item(61, 11).Value = '20/10/2023 02:12'
$ws.Cells.Item(61, 12).Value = 1.42
$ws.Cells.Item(61, 13).Value = '21/10/2023 14:58'
$ws.Cells.Item(61, 14).Value = 4.78
$ws.Cells.Item(61, 15).Value = '20/10/2023 02:12'
$ws.Cells.Item(61, 16).Value = 4.4
$ws.Cells.Item(61, 17).Value = '21/10/2023 14:58'
$ws.Cells.Item(61, 18).Value = 7.52
$ws.Cells.Item(61, 19).Value = '20/10/2023 02:12'
$ws.Cells.Item(61, 20).Value = 7.13
$ws.Cells.Item(61, 21).Value = '21/10/2023 14:58'
$ws.Cells.Item(61, 22).Value = 'https://www.betexplorer.com/football/montenegro/prva-crnogorska-liga/sutjeska-jedinstvo/MBIq7ThK/'

# Write row 63
$ws.Cells.Item(63, 6).Value = 'Arsenal Tivat'
$ws.Cells.Item(63, 7).Value = 2
$ws.Cells.Item(63, 8).Value = 'Mornar Bar'
$ws.Cells.Item(63, 9).Value = 2
$ws.Cells.Item(63, 10).Value = 2.29
$ws.Cells.Item(63, 11).Value = '20/10/2023 02:12'
$ws.Cells.Item(63, 12).Value = 2.72
$ws.Cells.Item(63, 13).Value = '21/10/2023 14:43'
$ws.Cells.Item(63, 14).Value = 2.75
$ws.Cells.Item(63, 15).Value = '20/10/2023 02:12'
$ws.Cells.Item(63, 16).Value = 2.49
$ws.Cells.Item(63, 17).Value = '21/10/2023 14:43'
$ws.Cells.Item(63, 18).Value = 3.2
$ws.Cells.Item(63, 19).Value = '20/10/2023 02:12'
$ws.Cells.Item(63, 20).Value = 3.32
$ws.Cells.Item(63, 21).Value = '21/10/2023 14:43'
$ws.Cells.Item(63, 22).Value = 'https://www.betexplorer.com/football/montenegro/prva-crnogorska-liga/arsenal-tivat-mornar-bar/rLJu89wE/'

# Write row 52
$ws.Cells.Item(52, 6).Value = 'Sutjeska'
$ws.Cells.Item(52, 7).Value = 0
$ws.Cells.Item(52, 8).Value = 'Mornar Bar'
$ws.Cells.Item(52, 9).Value = 1
$ws.Cells.Item(52, 10).Value = 1.41
$ws.Cells.Item(52, 11).Value = '01/10/2023 12:43'
$ws.Cells.Item(52, 12).Value = 1.51
$ws.Cells.Item(52, 13).Value = '01/10/2023 17:46'
$ws.Cells.Item(52, 14).Value = 4.18
$ws.Cells.Item(52, 15).Value = '01/10/2023 12:43'
$ws.Cells.Item(52, 16).Value = 3.88
$ws.Cells.Item(52, 17).Value = '01/10/2023 17:46'
$ws.Cells.Item(52, 18).Value = 7.18
$ws.Cells.Item(52, 19).Value = '01/10/2023 12:43'
$ws.Cells.Item(52, 20).Value = 6.55
$ws.Cells.Item(52, 21).Value = '01/10/2023 17:46'
$ws.Cells.Item(52, 22).Value = 'https://www.betexplorer.com/football/montenegro/prva-crnogorska-liga/sutjeska-mornar-bar/Gp0beWg8/'

# Write row 53
$ws.Cells.Item(53, 6).Value = 'Decic'
$ws.Cells.Item(53, 7).Value = 2
$ws.Cells.Item(53, 8).Value = 'Buducnost'
$ws.Cells.Item(53, 9).Value = 0
$ws.Cells.Item(53, 10).Value = 3.11
$ws.Cells.Item(53, 11).Value = '30/09/2023 12:43'
$ws.Cells.Item(53, 12).Value = 2.82
$ws.Cells.Item(53, 13).Value = '01/10/2023 17:39'
$ws.Cells.Item(53, 14).Value = 2.92
$ws.Cells.Item(53, 15).Value = '30/09/2023 12:43'
$ws.Cells.Item(53, 16).Value = 2.79
$ws.Cells.Item(53, 17).Value = '01/10/2023 17:39'
$ws.Cells.Item(53, 18).Value = 2.22
$ws.Cells.Item(53, 19).Value = '30/09/2023 12:43'
$ws.Cells.Item(53, 20).Value = 2.79
$ws.Cells.Item(53, 21).Value = '01/10/2023 17:39'
$ws.Cells.Item(53, 22).Value = 'https://www.betexplorer.com/football/montenegro/prva-crnogorska-liga/decic-buducnost/212AhUwR/'

# Write row 54
$ws.Cells.Item(54, 6).Value = 'Jezero'
$ws.Cells.Item(54, 7).Value = 1
$ws.Cells.Item(54, 8).Value = 'Arsenal Tivat'
$ws.Cells.Item(54, 9).Value = 1
$ws.Cells.Item(54, 10).Value = 2.08
$ws.Cells.Item(54, 11).Value = '30/09/2023 05:12'
$ws.Cells.Item(54, 12).Value = 2.41
$ws.Cells.Item(54, 13).Value = '01/10/2023 17:50'
$ws.Cells.Item(54, 14).Value = 2.84
$ws.Cells.Item(54, 15).Value = '30/09/2023 05:12'
$ws.Cells.Item(54, 16).Value = 2.84
$ws.Cells.Item(54, 17).Value = '01/10/2023 17:50'
$ws.Cells.Item(54, 18).Value = 3.55
$ws.Cells.Item(54, 19).Value = '30/09/2023 05:12'
$ws.Cells.Item(54, 20).Value = 3.3
$ws.Cells.Item(54, 21).Value = '01/10/2023 17:50'
$ws.Cells.Item(54, 22).Value = 'https://www.betexplorer.com/football/montenegro/prva-crnogorska-liga/jezero-arsenal-tivat/fya2fj9E/'

# Prepare formatting for new rows by copying from row 65 (columns A and E carry special styles)
$ws.Range("A65").Copy($ws.Range("A66"))
$ws.Range("E65").Copy($ws.Range("E66"))
$ws.Range("A65").Copy($ws.Range("A67"))
$ws.Range("E65").Copy($ws.Range("E67"))
$ws.Range("A65").Copy($ws.Range("A68"))
$ws.Range("E65").Copy($ws.Range("E68"))
$ws.Range("A65").Copy($ws.Range("A69"))
$ws.Range("E65").Copy($ws.Range("E69"))
$ws.Range("A65").Copy($ws.Range("A70"))
$ws.Range("E65").Copy($ws.Range("E70"))

# Write new row 66 (Indice=65)
$ws.Cells.Item(66, 1).Value = 65
$ws.Cells.Item(66, 2).Value = 'montenegro'
$ws.Cells.Item(66, 3).Value = 'prva-crnogorska-liga'
$ws.Cells.Item(66, 4).Value = '2023-2024'
$ws.Cells.Item(66, 5).Value = 45227.625
$ws.Cells.Item(66, 6).Value = 'Jedinstvo'
$ws.Cells.Item(66, 7).Value = 3
$ws.Cells.Item(66, 8).Value = 'Arsenal Tivat'
$ws.Cells.Item(66, 9).Value = 2
$ws.Cells.Item(66, 10).Value = 2.48
$ws.Cells.Item(66, 11).Value = '27/10/2023 03:12'
$ws.Cells.Item(66, 12).Value = 2.95
$ws.Cells.Item(66, 13).Value = '28/10/2023 14:00'
$ws.Cells.Item(66, 14).Value = 2.84
$ws.Cells.Item(66, 15).Value = '27/10/2023 03:12'
$ws.Cells.Item(66, 16).Value = 2.9
$ws.Cells.Item(66, 17).Value = '28/10/2023 14:00'
$ws.Cells.Item(66, 18).Value = 2.79
$ws.Cells.Item(66, 19).Value = '27/10/2023 03:12'
$ws.Cells.Item(66, 20).Value = 2.58
$ws.Cells.Item(66, 21).Value = '28/10/2023 14:00'
$ws.Cells.Item(66, 22).Value = 'https://www.betexplorer.com/football/montenegro/prva-crnogorska-liga/jedinstvo-arsenal-tivat/2q9vd3pJ/'

# Write new row 67 (Indice=66)
$ws.Cells.Item(67, 1).Value = 66
$ws.Cells.Item(67, 2).Value = 'montenegro'
$ws.Cells.Item(67, 3).Value = 'prva-crnogorska-liga'
$ws.Cells.Item(67, 4).Value = '2023-2024'
$ws.Cells.Item(67, 5).Value = 45227.66666666666
$ws.Cells.Item(67, 6).Value = 'Petrovac'
$ws.Cells.Item(67, 7).Value = 0
$ws.Cells.Item(67, 8).Value = 'Jezero'
$ws.Cells.Item(67, 9).Value = 0
$ws.Cells.Item(67, 10).Value = 2.03
$ws.Cells.Item(67, 11).Value = '27/10/2023 04:12'
$ws.Cells.Item(67, 12).Value = 2.39
$ws.Cells.Item(67, 13).Value = '28/10/2023 15:45'
$ws.Cells.Item(67, 14).Value = 2.98
$ws.Cells.Item(67, 15).Value = '27/10/2023 04:12'
$ws.Cells.Item(67, 16).Value = 2.84
$ws.Cells.Item(67, 17).Value = '28/10/2023 15:45'
$ws.Cells.Item(67, 18).Value = 3.49
$ws.Cells.Item(67, 19).Value = '27/10/2023 04:12'
$ws.Cells.Item(67, 20).Value = 3.34
$ws.Cells.Item(67, 21).Value = '28/10/2023 15:45'
$ws.Cells.Item(67, 22).Value = 'https://www.betexplorer.com/football/montenegro/prva-crnogorska-liga/petrovac-jezero/fJ5WcPF6/'

# Write new row 68 (Indice=67)
$ws.Cells.Item(68, 1).Value = 67
$ws.Cells.Item(68, 2).Value = 'montenegro'
$ws.Cells.Item(68, 3).Value = 'prva-crnogorska-liga'
$ws.Cells.Item(68, 4).Value = '2023-2024'
$ws.Cells.Item(68, 5).Value = 45227.70833333334
$ws.Cells.Item(68, 6).Value = 'Mornar Bar'
$ws.Cells.Item(68, 7).Value = 0
$ws.Cells.Item(68, 8).Value = 'Decic'
$ws.Cells.Item(68, 9).Value = 1
$ws.Cells.Item(68, 10).Value = 3.11
$ws.Cells.Item(68, 11).Value = '27/10/2023 05:12'
$ws.Cells.Item(68, 12).Value = 4.89
$ws.Cells.Item(68, 13).Value = '28/10/2023 16:07'
$ws.Cells.Item(68, 14).Value = 2.8
$ws.Cells.Item(68, 15).Value = '27/10/2023 05:12'
$ws.Cells.Item(68, 16).Value = 3.5
$ws.Cells.Item(68, 17).Value = '28/10/2023 16:59'
$ws.Cells.Item(68, 18).Value = 2.3
$ws.Cells.Item(68, 19).Value = '27/10/2023 05:12'
$ws.Cells.Item(68, 20).Value = 1.65
$ws.Cells.Item(68, 21).Value = '28/10/2023 16:29'
$ws.Cells.Item(68, 22).Value = 'https://www.betexplorer.com/football/montenegro/prva-crnogorska-liga/mornar-bar-decic/Y3HNaoof/'

# Write new row 69 (Indice=68)
$ws.Cells.Item(69, 1).Value = 68
$ws.Cells.Item(69, 2).Value = 'montenegro'
$ws.Cells.Item(69, 3).Value = 'prva-crnogorska-liga'
$ws.Cells.Item(69, 4).Value = '2023-2024'
$ws.Cells.Item(69, 5).Value = 45227.70833333334
$ws.Cells.Item(69, 6).Value = 'Mladost DG'
$ws.Cells.Item(69, 7).Value = 0
$ws.Cells.Item(69, 8).Value = 'Rudar'
$ws.Cells.Item(69, 9).Value = 1
$ws.Cells.Item(69, 10).Value = 2.17
$ws.Cells.Item(69, 11).Value = '27/10/2023 05:12'
$ws.Cells.Item(69, 12).Value = 1.53
$ws.Cells.Item(69, 13).Value = '28/10/2023 16:54'
$ws.Cells.Item(69, 14).Value = 2.98
$ws.Cells.Item(69, 15).Value = '27/10/2023 05:12'
$ws.Cells.Item(69, 16).Value = 3.93
$ws.Cells.Item(69, 17).Value = '28/10/2023 16:55'
$ws.Cells.Item(69, 18).Value = 3.15
$ws.Cells.Item(69, 19).Value = '27/10/2023 05:12'
$ws.Cells.Item(69, 20).Value = 5.06
$ws.Cells.Item(69, 21).Value = '28/10/2023 16:55'
$ws.Cells.Item(69, 22).Value = 'https://www.betexplorer.com/football/montenegro/prva-crnogorska-liga/mladost-dg-rudar/GS6Sb500/'

# Write new row 70 (Indice=69)
$ws.Cells.Item(70, 1).Value = 69
$ws.Cells.Item(70, 2).Value = 'montenegro'
$ws.Cells.Item(70, 3).Value = 'prva-crnogorska-liga'
$ws.Cells.Item(70, 4).Value = '2023-2024'
$ws.Cells.Item(70, 5).Value = 45227.77083333334
$ws.Cells.Item(70, 6).Value = 'Buducnost'
$ws.Cells.Item(70, 7).Value = 1
$ws.Cells.Item(70, 8).Value = 'Sutjeska'
$ws.Cells.Item(70, 9).Value = 0
$ws.Cells.Item(70, 10).Value = 2.04
$ws.Cells.Item(70, 11).Value = '27/10/2023 06:42'
$ws.Cells.Item(70, 12).Value = 2.13
$ws.Cells.Item(70, 13).Value = '28/10/2023 18:29'
$ws.Cells.Item(70, 14).Value = 3.13
$ws.Cells.Item(70, 15).Value = '27/10/2023 06:42'
$ws.Cells.Item(70, 16).Value = 3.17
$ws.Cells.Item(70, 17).Value = '28/10/2023 18:24'
$ws.Cells.Item(70, 18).Value = 3.29
$ws.Cells.Item(70, 19).Value = '27/10/2023 06:42'
$ws.Cells.Item(70, 20).Value = 3.5
$ws.Cells.Item(70, 21).Value = '28/10/2023 18:29'
$ws.Cells.Item(70, 22).Value = 'https://www.betexplorer.com/football/montenegro/prva-crnogorska-liga/buducnost-sutjeska/OhAzcqVC/'
